$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Regenerated s_vals data (filtered save games) for rows 2-20, columns B:E and G
$ws.Cells.Item(2, 2).Value = 0.04172184405617529
$ws.Cells.Item(2, 3).Value = 0.04103571897497393
$ws.Cells.Item(2, 4).Value = 0.7210945179870265
$ws.Cells.Item(2, 5).Value = 0.5333859586016987
$ws.Cells.Item(2, 7).Value = 1.337238039619874
$ws.Cells.Item(3, 2).Value = 3.272327238179451
$ws.Cells.Item(3, 3).Value = 1.626987699542094
$ws.Cells.Item(3, 4).Value = 0.1496068669990043
$ws.Cells.Item(3, 5).Value = 0.5333859586016987
$ws.Cells.Item(3, 7).Value = 5.582307763322248
$ws.Cells.Item(4, 2).Value = 0.04172184405617529
$ws.Cells.Item(4, 3).Value = 0.3048912486333797
$ws.Cells.Item(4, 4).Value = 0.7210945179870265
$ws.Cells.Item(4, 5).Value = 0.5333859586016987
$ws.Cells.Item(4, 7).Value = 1.60109356927828
$ws.Cells.Item(5, 2).Value = 0.6545652718822623
$ws.Cells.Item(5, 3).Value = 1.626987699542094
$ws.Cells.Item(5, 4).Value = 0.7210945179870265
$ws.Cells.Item(5, 5).Value = 0.5333859586016987
$ws.Cells.Item(5, 7).Value = 3.536033448013082
$ws.Cells.Item(6, 2).Value = 3.272327238179451
$ws.Cells.Item(6, 3).Value = 1.626987699542094
$ws.Cells.Item(6, 4).Value = 3.223369029078222
$ws.Cells.Item(6, 5).Value = 0.5333859586016987
$ws.Cells.Item(6, 7).Value = 8.656069925401464
$ws.Cells.Item(7, 2).Value = 3.272327238179451
$ws.Cells.Item(7, 3).Value = 1.626987699542094
$ws.Cells.Item(7, 4).Value = 0.1496068669990043
$ws.Cells.Item(7, 5).Value = 0.5333859586016987
$ws.Cells.Item(7, 7).Value = 5.582307763322248
$ws.Cells.Item(8, 2).Value = 1.445647641019636
$ws.Cells.Item(8, 3).Value = 1.626987699542094
$ws.Cells.Item(8, 4).Value = 3.223369029078222
$ws.Cells.Item(8, 5).Value = 0.5333859586016987
$ws.Cells.Item(8, 7).Value = 6.82939032824165
$ws.Cells.Item(9, 2).Value = 3.272327238179451
$ws.Cells.Item(9, 3).Value = 1.626987699542094
$ws.Cells.Item(9, 4).Value = 0.1496068669990043
$ws.Cells.Item(9, 5).Value = 0.5333859586016987
$ws.Cells.Item(9, 7).Value = 5.582307763322248
$ws.Cells.Item(10, 2).Value = 1.445647641019636
$ws.Cells.Item(10, 3).Value = 1.626987699542094
$ws.Cells.Item(10, 4).Value = 0.7210945179870265
$ws.Cells.Item(10, 5).Value = 0.5333859586016987
$ws.Cells.Item(10, 7).Value = 4.327115817150455
$ws.Cells.Item(11, 2).Value = 1.445647641019636
$ws.Cells.Item(11, 3).Value = 1.626987699542094
$ws.Cells.Item(11, 4).Value = 0.7210945179870265
$ws.Cells.Item(11, 5).Value = 0.5333859586016987
$ws.Cells.Item(11, 7).Value = 4.327115817150455
$ws.Cells.Item(12, 2).Value = 1.445647641019636
$ws.Cells.Item(12, 3).Value = 1.626987699542094
$ws.Cells.Item(12, 4).Value = 0.1496068669990043
$ws.Cells.Item(12, 5).Value = 0.5333859586016987
$ws.Cells.Item(12, 7).Value = 3.755628166162433
$ws.Cells.Item(13, 2).Value = 0.6545652718822623
$ws.Cells.Item(13, 3).Value = 1.626987699542094
$ws.Cells.Item(13, 4).Value = 3.223369029078222
$ws.Cells.Item(13, 5).Value = 0.5333859586016987
$ws.Cells.Item(13, 7).Value = 6.038307959104277
$ws.Cells.Item(14, 2).Value = 3.272327238179451
$ws.Cells.Item(14, 3).Value = 1.626987699542094
$ws.Cells.Item(14, 4).Value = 18.71679738969934
$ws.Cells.Item(14, 5).Value = 13.86384647080068
$ws.Cells.Item(14, 7).Value = 37.47995879822157
$ws.Cells.Item(15, 2).Value = 3.272327238179451
$ws.Cells.Item(15, 3).Value = 1.626987699542094
$ws.Cells.Item(15, 4).Value = 0.7210945179870265
$ws.Cells.Item(15, 5).Value = 0.5333859586016987
$ws.Cells.Item(15, 7).Value = 6.15379541431027
$ws.Cells.Item(16, 2).Value = 3.272327238179451
$ws.Cells.Item(16, 3).Value = 1.626987699542094
$ws.Cells.Item(16, 4).Value = 0.1496068669990043
$ws.Cells.Item(16, 5).Value = 0.5333859586016987
$ws.Cells.Item(16, 7).Value = 5.582307763322248
$ws.Cells.Item(17, 2).Value = 1.445647641019636
$ws.Cells.Item(17, 3).Value = 0.3048912486333797
$ws.Cells.Item(17, 4).Value = 0.1496068669990043
$ws.Cells.Item(17, 5).Value = 0.5333859586016987
$ws.Cells.Item(17, 7).Value = 2.433531715253719
$ws.Cells.Item(18, 2).Value = 1.445647641019636
$ws.Cells.Item(18, 3).Value = 1.626987699542094
$ws.Cells.Item(18, 4).Value = 0.7210945179870265
$ws.Cells.Item(18, 5).Value = 0.5333859586016987
$ws.Cells.Item(18, 7).Value = 4.327115817150455
$ws.Cells.Item(19, 2).Value = 3.272327238179451
$ws.Cells.Item(19, 3).Value = 1.626987699542094
$ws.Cells.Item(19, 4).Value = 0.7210945179870265
$ws.Cells.Item(19, 5).Value = 0.5333859586016987
$ws.Cells.Item(19, 7).Value = 6.15379541431027
$ws.Cells.Item(20, 2).Value = 3.272327238179451
$ws.Cells.Item(20, 3).Value = 1.626987699542094
$ws.Cells.Item(20, 4).Value = 0.7210945179870265
$ws.Cells.Item(20, 5).Value = 0.5333859586016987
$ws.Cells.Item(20, 7).Value = 6.15379541431027

$wb.Save()
